$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("H2").Value = 4.33
$ws.Range("I2").Value = 6.25
$ws.Range("AH2").Value = 17
$ws.Range("AJ2").Value = 19
$ws.Range("AO2").Value = 7.5
$ws.Range("AW2").Value = 7.5
$ws.Range("AZ2").Value = 101

# Row 4 updates
$ws.Range("M4").Value = 1.11
$ws.Range("N4").Value = 6.5
